$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 15 de Agosto de 2020 a las 15:46"

# Country name changes caused by re-sorting (ranking) of the table
$ws.Range("A64").Value2 = "Kenia"
$ws.Range("A65").Value2 = "Moldavia"
$ws.Range("A88").Value2 = "Zambia"
$ws.Range("A89").Value2 = "Malasia"
$ws.Range("A90").Value2 = "Paraguay"
$ws.Range("A128").Value2 = "Eslovenia"
$ws.Range("A129").Value2 = "Lituania"
$ws.Range("A213").Value2 = "Montserrat"
$ws.Range("A214").Value2 = "Islas Malvinas"

# Updated statistics for countries whose figures changed
$ws.Range("B6").Value2 = 2545062
$ws.Range("C6").Value2 = 19840
$ws.Range("D6").Value2 = 1824940
$ws.Range("E6").Value2 = 670799
$ws.Range("G6").Value2 = 189
$ws.Range("H6").Value2 = 49323
$ws.Range("B16").Value2 = 297315
$ws.Range("C16").Value2 = 1413
$ws.Range("D16").Value2 = 264487
$ws.Range("E16").Value2 = 29459
$ws.Range("G16").Value2 = 31
$ws.Range("H16").Value2 = 3369
$ws.Range("B24").Value2 = 172583
$ws.Range("C24").Value2 = 4293
$ws.Range("D24").Value2 = 122700
$ws.Range("E24").Value2 = 44098
$ws.Range("G24").Value2 = 76
$ws.Range("H24").Value2 = 5785
$ws.Range("B45").Value2 = 62495
$ws.Range("C45").Value2 = 655
$ws.Range("G45").Value2 = 2
$ws.Range("H45").Value2 = 6169
$ws.Range("B49").Value2 = 53981
$ws.Range("C49").Value2 = 198
$ws.Range("D49").Value2 = 39585
$ws.Range("E49").Value2 = 12621
$ws.Range("G49").Value2 = 3
$ws.Range("H49").Value2 = 1775
$ws.Range("B61").Value2 = 34107
$ws.Range("C61").Value2 = 89
$ws.Range("D61").Value2 = 31697
$ws.Range("E61").Value2 = 1904
$ws.Range("G61").Value2 = 2
$ws.Range("H61").Value2 = 506
$ws.Range("B64").Value2 = 29849
$ws.Range("C64").Value2 = 515
$ws.Range("D64").Value2 = 15298
$ws.Range("E64").Value2 = 14086
$ws.Range("H64").Value2 = 465
$ws.Range("B65").Value2 = 29483
$ws.Range("D65").Value2 = 20556
$ws.Range("E65").Value2 = 8043
$ws.Range("H65").Value2 = 884
$ws.Range("B66").Value2 = 29471
$ws.Range("C66").Value2 = 238
$ws.Range("E66").Value2 = 2361
$ws.Range("G66").Value2 = 5
$ws.Range("H66").Value2 = 670
$ws.Range("B78").Value2 = 15801
$ws.Range("C78").Value2 = 266
$ws.Range("D78").Value2 = 9619
$ws.Range("E78").Value2 = 5711
$ws.Range("G78").Value2 = 2
$ws.Range("H78").Value2 = 471
$ws.Range("B79").Value2 = 15483
$ws.Range("C79").Value2 = 104
$ws.Range("D79").Value2 = 13275
$ws.Range("E79").Value2 = 1587
$ws.Range("B86").Value2 = 9934
$ws.Range("C86").Value2 = 26
$ws.Range("E86").Value2 = 816
$ws.Range("B88").Value2 = 9186
$ws.Range("C88").Value2 = 165
$ws.Range("D88").Value2 = 8065
$ws.Range("E88").Value2 = 861
$ws.Range("G88").Value2 = 4
$ws.Range("H88").Value2 = 260
$ws.Range("B89").Value2 = 9175
$ws.Range("C89").Value2 = 26
$ws.Range("D89").Value2 = 8831
$ws.Range("E89").Value2 = 219
$ws.Range("H89").Value2 = 125
$ws.Range("B90").Value2 = 9022
$ws.Range("D90").Value2 = 5657
$ws.Range("E90").Value2 = 3257
$ws.Range("H90").Value2 = 108
$ws.Range("D111").Value2 = 3488
$ws.Range("E111").Value2 = 851
$ws.Range("G111").Value2 = 1
$ws.Range("H111").Value2 = 68
$ws.Range("B128").Value2 = 2401
$ws.Range("C128").Value2 = 32
$ws.Range("D128").Value2 = 2042
$ws.Range("E128").Value2 = 230
$ws.Range("H128").Value2 = 129
$ws.Range("B129").Value2 = 2386
$ws.Range("C129").Value2 = 34
$ws.Range("D129").Value2 = 1704
$ws.Range("E129").Value2 = 601
$ws.Range("H129").Value2 = 81
$ws.Range("D213").Value2 = 12
$ws.Range("H213").Value2 = 1
$ws.Range("D214").Value2 = 13
$ws.Range("H214").Value2 = 0
